$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 72.45
$ws.Range("I15").Value = 72.45
$ws.Range("K15").Value = 217.35
$ws.Range("M15").Value = -48.35000000000002

$ws.Range("H40").Value = 1047.9032
$ws.Range("I40").Value = 999.6
$ws.Range("J40").Value = 1057.1923
$ws.Range("K40").Value = 999.6
$ws.Range("L40").Value = 1057.1923
$ws.Range("M40").Value = -824.6
$ws.Range("N40").Value = -1407.1923

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H138").Value = 2780884.8
$ws.Range("I138").Value = 2262
$ws.Range("J138").Value = 5132027
$ws.Range("K138").Value = 6786
$ws.Range("L138").Value = 15396081
$ws.Range("M138").Value = -1646
$ws.Range("N138").Value = -15406361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

$ws.Range("H13").Value = 23200560
$ws.Range("I13").Value = 29000000
$ws.Range("J13").Value = 2800
$ws.Range("K13").Value = 29000000
$ws.Range("L13").Value = 2800
$ws.Range("M13").Value = -28999856
$ws.Range("N13").Value = -3088

$ws.Range("H32").Value = 30812.043
$ws.Range("I32").Value = 30904.035
$ws.Range("J32").Value = 30667.838
$ws.Range("K32").Value = 30904.035
$ws.Range("L32").Value = 30667.838
$ws.Range("M32").Value = -30617.035
$ws.Range("N32").Value = -31241.838

$ws.Range("H74").Value = 9865556
$ws.Range("I74").Value = 14537670
$ws.Range("J74").Value = 96590.91
$ws.Range("K74").Value = 14537670
$ws.Range("L74").Value = 96590.91
$ws.Range("M74").Value = -14536796
$ws.Range("N74").Value = -98338.91

$ws.Range("H77").Value = 9865556
$ws.Range("I77").Value = 14537670
$ws.Range("J77").Value = 96590.91
$ws.Range("K77").Value = 72688350
$ws.Range("L77").Value = 482954.55
$ws.Range("M77").Value = -72683982
$ws.Range("N77").Value = -491690.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 438.5
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 477
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 477
$ws.Range("M12").Value = -232
$ws.Range("N12").Value = -813

$ws.Range("H14").Value = 3000
$ws.Range("J14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("N14").Value = -3344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 6174.6665
$ws.Range("I15").Value = 3000
$ws.Range("J15").Value = 7762
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 7762
$ws.Range("M15").Value = -2830
$ws.Range("N15").Value = -8102

$ws.Range("H31").Value = 71555.09
$ws.Range("I31").Value = 74669.07000000001
$ws.Range("J31").Value = 69133.11
$ws.Range("K31").Value = 74669.07000000001
$ws.Range("L31").Value = 69133.11
$ws.Range("M31").Value = -74374.07000000001
$ws.Range("N31").Value = -69723.11

$ws.Range("H34").Value = 71555.09
$ws.Range("I34").Value = 74669.07000000001
$ws.Range("J34").Value = 69133.11
$ws.Range("K34").Value = 74669.07000000001
$ws.Range("L34").Value = 69133.11
$ws.Range("M34").Value = -74467.07000000001
$ws.Range("N34").Value = -69537.11

$ws.Range("H41").Value = 16258.333
$ws.Range("I41").Value = 6550
$ws.Range("J41").Value = 18200
$ws.Range("K41").Value = 6550
$ws.Range("L41").Value = 18200
$ws.Range("M41").Value = -6122
$ws.Range("N41").Value = -19056

$ws.Range("H50").Value = 23819
$ws.Range("J50").Value = 23819
$ws.Range("L50").Value = 23819
$ws.Range("N50").Value = -25069

$ws.Range("H51").Value = 25479.2
$ws.Range("J51").Value = 25479.2
$ws.Range("L51").Value = 25479.2
$ws.Range("N51").Value = -26951.2

$ws.Range("H61").Value = 25479.2
$ws.Range("J61").Value = 25479.2
$ws.Range("L61").Value = 25479.2
$ws.Range("N61").Value = -26175.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 797.875
$ws.Range("I34").Value = 341.66666
$ws.Range("J34").Value = 949.94446
$ws.Range("K34").Value = 1024.99998
$ws.Range("L34").Value = 2849.83338
$ws.Range("M34").Value = -940.9999800000001
$ws.Range("N34").Value = -3017.83338

$ws.Range("H39").Value = 7000
$ws.Range("J39").Value = 7000
$ws.Range("L39").Value = 21000
$ws.Range("N39").Value = -21588

$ws.Range("H55").Value = 2214.2856
$ws.Range("J55").Value = 2500
$ws.Range("L55").Value = 7500
$ws.Range("N55").Value = -7854

$ws.Range("H75").Value = 3100.8333
$ws.Range("I75").Value = 1300
$ws.Range("J75").Value = 4001.25
$ws.Range("K75").Value = 3900
$ws.Range("L75").Value = 12003.75
$ws.Range("M75").Value = -2902
$ws.Range("N75").Value = -13999.75

$ws.Range("H78").Value = 3100.8333
$ws.Range("I78").Value = 1300
$ws.Range("J78").Value = 4001.25
$ws.Range("K78").Value = 11700
$ws.Range("L78").Value = 36011.25
$ws.Range("M78").Value = -6708
$ws.Range("N78").Value = -45995.25

$ws.Range("H137").Value = 2883.85
$ws.Range("I137").Value = 1198.2858
$ws.Range("J137").Value = 3791.4614
$ws.Range("K137").Value = 3594.8574
$ws.Range("L137").Value = 11374.3842
$ws.Range("M137").Value = 1505.1426
$ws.Range("N137").Value = -21574.3842

$ws.Range("H140").Value = 3160
$ws.Range("I140").Value = 3166.6667
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 9500.000100000001
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -4320.000100000001
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5000000
$ws.Range("I10").Value = 5000000
$ws.Range("K10").Value = 5000000
$ws.Range("M10").Value = -4999831

$ws.Range("H46").Value = 16433.334
$ws.Range("I46").Value = 4650
$ws.Range("J46").Value = 40000
$ws.Range("K46").Value = 4650
$ws.Range("L46").Value = 40000
$ws.Range("M46").Value = -4494
$ws.Range("N46").Value = -40312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 10000003
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H55").Value = 217.60869
$ws.Range("I55").Value = 235.85
$ws.Range("J55").Value = 96
$ws.Range("K55").Value = 235.85
$ws.Range("L55").Value = 96
$ws.Range("M55").Value = -62.84999999999999
$ws.Range("N55").Value = -442

$ws.Range("H110").Value = 8528.799999999999
$ws.Range("J110").Value = 8528.799999999999
$ws.Range("L110").Value = 8528.799999999999
$ws.Range("N110").Value = -16708.8

$ws.Range("H127").Value = 58536
$ws.Range("J127").Value = 58536
$ws.Range("L127").Value = 58536
$ws.Range("N127").Value = -68456

$ws.Range("H136").Value = 28879.908
$ws.Range("I136").Value = 19272.836
$ws.Range("J136").Value = 65506.875
$ws.Range("K136").Value = 57818.508
$ws.Range("L136").Value = 196520.625
$ws.Range("M136").Value = -55268.508
$ws.Range("N136").Value = -201620.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23750
$ws.Range("J64").Value = 23750
$ws.Range("L64").Value = 23750
$ws.Range("N64").Value = -24246

$ws.Range("H67").Value = 23750
$ws.Range("J67").Value = 23750
$ws.Range("L67").Value = 23750
$ws.Range("N67").Value = -25466

$ws.Range("H132").Value = 40682.41
$ws.Range("I132").Value = 32280.625
$ws.Range("J132").Value = 54832.79
$ws.Range("K132").Value = 96841.875
$ws.Range("L132").Value = 164498.37
$ws.Range("M132").Value = -94311.875
$ws.Range("N132").Value = -169558.37
